$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Mapping of row number -> new value for column G (K = strikeouts),
# regenerated from the K stat instead of the old running Strike# counter.
$gUpdates = @{
    2 = 1
    3 = 1
    4 = 1
    5 = 0
    6 = 1
    7 = 1
    8 = 1
    9 = 1
    10 = 2
    11 = 2
    12 = 1
    13 = 2
    14 = 1
    15 = 1
    17 = 0
    18 = 1
    19 = 0
    20 = 1
    21 = 0
    22 = 0
    23 = 1
    24 = 0
    25 = 1
    26 = 1
    27 = 1
    28 = 1
    29 = 0
    30 = 1
    31 = 0
    32 = 2
    33 = 2
    34 = 1
    35 = 3
    36 = 2
    37 = 0
    38 = 1
    39 = 0
    40 = 0
    41 = 0
    42 = 0
    43 = 1
    44 = 0
    45 = 2
    46 = 1
    47 = 1
    48 = 1
    50 = 1
    51 = 0
    52 = 1
    53 = 1
    54 = 2
    55 = 0
    56 = 3
    57 = 0
    58 = 1
    59 = 1
    60 = 0
    61 = 0
    62 = 2
    63 = 1
    64 = 1
    65 = 0
    66 = 1
    67 = 1
    68 = 1
    69 = 2
    70 = 1
    71 = 2
    72 = 1
    73 = 1
    74 = 0
    75 = 1
    76 = 3
    77 = 0
    78 = 1
    80 = 1
    81 = 1
    82 = 1
    83 = 0
    84 = 1
}

foreach ($row in $gUpdates.Keys) {
    $ws.Cells.Item($row, 7).Value = $gUpdates[$row]
}
